# Adds wishlist functionality test sheets:
#  - Renames existing sheets to the ShoppingCart-qualified names
#  - Duplicates them to create matching Wishlist sheets
#  - Updates the TestCaseName column text on the new sheets
#  - Restores per-sheet selections to match the target workbook
#
# NOTE: worksheet handles are re-fetched by name right before each use.
# Copying/inserting sheets shifts collection indices, and stale handles
# obtained before such a shift can silently resolve to the wrong sheet.

$wb = $excel.ActiveWorkbook

# --- Rename the two existing sheets ---------------------------------------
$wb.Worksheets.Item("AddProduct").Name = "AddProductShoppingCart"
$wb.Worksheets.Item("DeleteProduct").Name = "DeleteProductFromShoppingCart"

# --- Duplicate "AddProductShoppingCart" -> "AddProductWishlist" ------------
$addCart = $wb.Worksheets.Item("AddProductShoppingCart")
$addCart.Copy($null, $addCart)
$wb.Worksheets.Item("AddProductShoppingCart (2)").Name = "AddProductWishlist"

$addWish = $wb.Worksheets.Item("AddProductWishlist")
$addWish.Range("C2").Value = "Add Hair Care Product to Wishlist"
$addWish.Range("C3").Value = "Add Apparel & accessories Product to Wishlist"
$addWish.Range("C4").Value = "Add Makeup Product to Wishlist"
$addWish.Range("C5").Value = "Add Skincare Product to Wishlist"
$addWish.Range("C6").Value = "Add Mens Product to Wishlist"
$addWish.Range("C7").Value = "Add Books Product to Wishlist"
$addWish.Range("C8").Value = "Add Fragrance Product to Wishlist"

# --- Duplicate "DeleteProductFromShoppingCart" -> "DeleteProductFromWishlist"
$delCart = $wb.Worksheets.Item("DeleteProductFromShoppingCart")
$delCart.Copy($null, $delCart)
$wb.Worksheets.Item("DeleteProductFromShoppingCart (2)").Name = "DeleteProductFromWishlist"

$delWish = $wb.Worksheets.Item("DeleteProductFromWishlist")
$delWish.Range("C2").Value = "Delete Pantene Pro-V Product from Wish list"
$delWish.Range("C3").Value = "Delete Shaving cream Product from Wish list"

# --- Restore selections on each sheet --------------------------------------
$wb.Worksheets.Item("AddProductShoppingCart").Range("C15").Select()
$wb.Worksheets.Item("DeleteProductFromShoppingCart").Range("C10").Select()
$wb.Worksheets.Item("DeleteProductFromWishlist").Range("G13").Select()

# --- Make the Wishlist-add sheet the active tab, with its own selection ----
$wb.Worksheets.Item("AddProductWishlist").Activate()
$wb.Worksheets.Item("AddProductWishlist").Range("E15").Select()
